$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New "Fix Groundcheck" task description, reused later for the new row 9.
$newTaskText = "Fix Groundcheck issues in player controller that causes groundcheck to sometimes fail groundcheck"

# Mirror column H (the "Past" week's Friday column) into a brand new
# column I - i.e. a fresh "week" section ("Make a new scene").
$ws.Range("I2").Value2 = $ws.Range("H2").Value2
$ws.Range("I3").Value2 = $ws.Range("H3").Value2
$ws.Range("I4").Value2 = $ws.Range("H4").Value2
$ws.Range("I5").Value2 = $ws.Range("H5").Value2
$ws.Range("I7").Value2 = $ws.Range("H7").Value2
$ws.Range("I8").Value2 = $ws.Range("H8").Value2

# Add the brand-new task as a new row in that same column.
$ws.Range("I9").Value2 = $newTaskText

# Widen column I to fit the new (long) content, matching the other
# "bestFit" columns on the sheet.
$ws.Columns.Item(9).ColumnWidth = 90.2

# Move the view over to the new column and select the new cell, like the
# author would have after typing the new task in.
$ws.Range("G1").Select()
$ws.Range("I9").Select()
